$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G16").Value = 80
$ws.Range("G18").Value = 89
$ws.Range("G23").Value = 124
$ws.Range("G26").Value = 261
$ws.Range("G27").Value = 327
$ws.Range("G28").Value = 297
$ws.Range("G29").Value = 345
$ws.Range("G30").Value = 368
$ws.Range("G31").Value = 440
$ws.Range("G32").Value = 411
$ws.Range("G33").Value = 403
$ws.Range("G34").Value = 519
$ws.Range("G36").Value = 451
$ws.Range("G37").Value = 531
$ws.Range("G38").Value = 428
$ws.Range("G39").Value = 624
$ws.Range("G40").Value = 667
$ws.Range("G41").Value = 535
$ws.Range("G42").Value = 651
$ws.Range("G43").Value = 687
$ws.Range("G44").Value = 690
$ws.Range("G46").Value = 736
$ws.Range("G47").Value = 899
$ws.Range("G48").Value = 801
